$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 30 (inclusive), shifting remaining rows up.
$ws.Rows("2:30").Delete()

$ws.Range("C11").Select()
